$d = $word.ActiveDocument

# 1. Update the Author line: "Jessie Zhao" -> "Kevin Arne" and
#    "Prototyping Labs Technical Assistance" -> "Prototyping Lab Supervisor"
$d.Content.Find.Execute("Jessie Zhao", $false, $false, $false, $false, $false, $true, 1, $false, "Kevin Arne", 2) | Out-Null
$d.Content.Find.Execute("Prototyping Labs Technical Assistance", $false, $false, $false, $false, $false, $true, 1, $false, "Prototyping Lab Supervisor", 2) | Out-Null

# 2. Clean up "(if applicable)" so it is a single run (removing the
#    gramStart/gramEnd proof-reading split around "if").
$d.Content.Find.Execute("(if applicable)", $false, $false, $false, $false, $false, $true, 1, $false, "(if applicable)", 2) | Out-Null

# 3. Header: "Prototyping Labs at GIX" -> "Prototyping Lab at GIX"
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.Find.Execute("Prototyping Labs at GIX", $false, $false, $false, $false, $false, $true, 1, $false, "Prototyping Lab at GIX", 2) | Out-Null

Write-Host "Edit complete"
